$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.596.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = "'2.343.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.05%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'558.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.31%  '
$ws.Range("D6").Value = "'132.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.56%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = "'0.579"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.91%  '
$ws.Range("E9").Value = '  -1.31%  '
$ws.Range("E10").Value = '  -2.05%  '
$ws.Range("E11").Value = '  +0.84%  '
$ws.Range("D12").Value = "'0.339"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.14%  '
$ws.Range("D13").Value = "'23.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.63%  '
$ws.Range("D14").Value = "'2.764.46"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.95%  '
$ws.Range("D15").Value = "'59.588.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.61%  '
$ws.Range("E16").Value = '  -0.76%  '
$ws.Range("D17").Value = "'2.350.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.17%  '
$ws.Range("D18").Value = "'10.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.10%  '
$ws.Range("E19").Value = '  +0.20%  '
$ws.Range("D20").Value = "'318.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'6.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.75%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = "'63.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.87%  '
$ws.Range("E24").Value = '  -2.96%  '
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").Value = "'8.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = '  -1.94%  '
$ws.Range("E28").Value = '  +1.84%  '
$ws.Range("D29").Value = "'171.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.44%  '
$ws.Range("D30").Value = "'0.0₃0749"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.62%  '
$ws.Range("D31").Value = "'5.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.93%  '
$ws.Range("E32").Value = '  +6.73%  '
$ws.Range("D33").Value = "'0.398"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.71%  '
$ws.Range("D34").Value = "'18.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.05%  '
$ws.Range("E36").Value = '  -1.28%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("E38").Value = '  -3.04%  '
$ws.Range("E39").Value = '  -1.96%  '
$ws.Range("D40").Value = "'314.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.21%  '
$ws.Range("D41").Value = "'38.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.34%  '
$ws.Range("D42").Value = "'144.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.19%  '
$ws.Range("E43").Value = '  -5.15%  '
$ws.Range("D44").Value = "'0.0958"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.47%  '
$ws.Range("E45").Value = '  -1.89%  '
$ws.Range("D46").Value = "'18.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.32%  '
$ws.Range("D47").Value = "'0.562"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.61%  '
$ws.Range("E48").Value = '  -3.04%  '
$ws.Range("D49").Value = "'11.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("E51").Value = '  -0.34%  '
